$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Locator Name"); the old column B
# ("Locator", e.g. "id=txtUsername") shifts to column C and is renamed
# "Locator value", keeping only the value portion (the "id=" prefix is
# split out into the new "Locator Name" column).
$ws.Columns("B").Insert()

# The inserted column doesn't inherit the old column B's width, so match
# column C's (old column B's) width onto the new column B.
$ws.Range("B1").ColumnWidth = $ws.Range("C1").ColumnWidth

# Header row
$ws.Range("A1").Value = "Test Step"
$ws.Range("B1").Value = "Locator Name"
$ws.Range("C1").Value = "Locator value"
$ws.Range("D1").Value = "Action"
$ws.Range("E1").Value = "Value"

# Row 2 - open browser
$ws.Range("A2").Value = "open browser"
$ws.Range("B2").Value = "NA"
$ws.Range("C2").Value = "NA"
$ws.Range("D2").Value = "open browser"
$ws.Range("E2").Value = "chrome"

# Row 3 - launch url
$ws.Range("A3").Value = "launch url"
$ws.Range("B3").Value = "NA"
$ws.Range("C3").Value = "NA"
$ws.Range("D3").Value = "launch url"
$ws.Range("E3").Value = "https://opensource-demo.orangehrmlive.com/"

# Row 4 - enter username
$ws.Range("A4").Value = "enter username"
$ws.Range("B4").Value = "id"
$ws.Range("C4").Value = "txtUsername"
$ws.Range("D4").Value = "sendkeys"
$ws.Range("E4").Value = "Admin"

# Row 5 - enter password
$ws.Range("A5").Value = "enter password"
$ws.Range("B5").Value = "id"
$ws.Range("C5").Value = "txtPassword"
$ws.Range("D5").Value = "sendkeys"
$ws.Range("E5").Value = "admin123"

# Row 6 - click on login button
$ws.Range("A6").Value = "click on login button"
$ws.Range("B6").Value = "id"
$ws.Range("C6").Value = "btnLogin"
$ws.Range("D6").Value = "click"
$ws.Range("E6").Value = "NA"

# Row 7 - close browser
$ws.Range("A7").Value = "close browser"
$ws.Range("B7").Value = "NA"
$ws.Range("C7").Value = "NA"
$ws.Range("D7").Value = "quit"
$ws.Range("E7").Value = "NA"

# The hyperlink that lived on D3 (launch-url row's Value cell) now lives on
# E3 since that column shifted right.
$ws.Range("D3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E3"), "https://opensource-demo.orangehrmlive.com/")

$ws.Range("E11").Select()
